$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their original text formatting instead of
# being auto-converted to numbers by Excel when the new values look numeric
# (e.g. "0.3930" must stay text, not become 0.393).
$priceCells = @("D2", "D3", "D5", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.542.30"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").Value = "1.693.11"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "315.44"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "0.3930"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "1.518"
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "52.78"
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("D12").Value = "0.08725"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "7.194"
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("D14").Value = "23.04"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "0.00001313"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "7.561"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "1.697.57"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "99.63"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "0.07051"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").Value = "19.57"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("D21").Value = "6.852"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "14.01"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "24.537.33"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").Value = "3.018"
$ws.Range("E25").Value = "  +7.18%  "
$ws.Range("D26").Value = "2.322"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "22.24"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").Value = "160.75"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "5.220"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "133.92"
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("D31").Value = "7.478"
$ws.Range("E31").Value = "  +11.06%  "
$ws.Range("D32").Value = "1.882.19"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "1.085"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").Value = "7.267"
$ws.Range("E34").Value = "  +10.02%  "
$ws.Range("D35").Value = "0.08506"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "11.25"
$ws.Range("E36").Value = "  +7.92%  "
$ws.Range("D37").Value = "1.949"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "0.2697"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "14.36"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "0.02742"
$ws.Range("E40").Value = "  +9.02%  "
$ws.Range("D41").Value = "0.09001"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").Value = "1.470"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").Value = "0.7607"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "0.7132"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "2.518"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "15.21"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "4.199"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "140.82"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "1.315"
$ws.Range("E50").Value = "  +6.04%  "
$ws.Range("D51").Value = "0.07986"
$ws.Range("E51").Value = "  +2.68%  "
